$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph via Find.
$rng = $d.Content
$found = $rng.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Ver no Jupiter...' paragraph to remove."
}

$fStart = $rng.Start

# Map that text position back to its paragraph index (Find's resulting Range
# doesn't reliably expose .Paragraphs/.Next/.Previous in this host, so walk
# the Paragraphs collection and compare ranges instead).
$idx = 0
$matchIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($fStart -ge $p.Range.Start -and $fStart -lt $p.Range.End) {
        $matchIdx = $idx
    }
}

if ($matchIdx -eq -1) {
    throw "Could not resolve the paragraph index for the found text."
}

# Remove three paragraphs in total: the blank separator paragraph right
# before "Ver no Jupiter...", the "Ver no Jupiter..." paragraph itself, and
# the footer/copyright paragraph right after it.
$startPara = $d.Paragraphs.Item($matchIdx - 1)
$endPara = $d.Paragraphs.Item($matchIdx + 1)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
